# Update to correlation results and fig caption
#
# Table S1 - Plasticity AIC (first worksheet): insert a new row for the
# "pstr" species using the full reef environment * pCO2 * temperature
# interaction model (re-run / footnoted model, marked with a trailing "7"
# footnote marker), pushing the existing pstr/past model rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank row at row 9 (pushes old rows 9-17 down to 10-18)
$ws.Rows.Item(9).Insert()

$ws.Range("A9").Value = "pstr"
$ws.Range("B9").Value = "reef environment * pCO2 * temperature + (1 | colony)7"
$ws.Range("C9").Value = "glmerMod"
$ws.Range("D9").Value = 105.1
$ws.Range("E9").Value = 116.4
$ws.Range("F9").Value = 0.3271
$ws.Range("G9").Value = 0.2594
$ws.Range("H9").Value = 0.0914015570605844
$ws.Range("I9").Value = 1.12817097438964
$ws.Range("J9").Value = 0.354967986837399
$ws.Range("K9").Value = 47.3

Write-Output "Inserted new pstr model row (row 9) in Table S1"
